# Generate Report for Handback
# Refresh the handoff/handback timestamps recorded on the Overview,
# zh-cn and de-de sheets of the handback status report.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Overview sheet: Latest HO Xliff Generate Date for 5faf0dda...md
$wsOverview.Range("G3").Value = "2016-08-23 14:54:53"

# zh-cn sheet: Correspond Handoff/Handback Datetime for 5faf0dda...md
$wsZhCn.Range("H3").Value = "2016-08-23 14:54:47"
$wsZhCn.Range("K3").Value = "2016-08-23 14:55:42"

# de-de sheet: Correspond Handoff/Handback Datetime for 5faf0dda...md
$wsDeDe.Range("H3").Value = "2016-08-23 14:54:53"
$wsDeDe.Range("K3").Value = "2016-08-23 14:55:50"
